# Updates cryptos list: refreshed Price (D) / Volume(1h) (E) figures for
# every coin row, and re-sorts two adjacent coin pairs whose ranking swapped
# (Bittensor <-> Binance-PegBSC-USD at rows 28/29, Stacks <-> Monero at rows 41/42)
# by rewriting their Coin (B) / Link (C) cells along with D/E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.018.34'
$ws.Range("E2").Value = '  -4.34%  '
$ws.Range("D3").Value = '2.449.85'
$ws.Range("E3").Value = '  -6.90%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''544.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.75%  '
$ws.Range("D6").Value = '''145.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.09%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''0.583'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.47%  '
$ws.Range("D9").Value = '2.450.91'
$ws.Range("E9").Value = '  -6.79%  '
$ws.Range("E10").Value = '  -10.81%  '
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").Value = '''5.40'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.07%  '
$ws.Range("D13").Value = '''0.349'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -9.36%  '
$ws.Range("D14").Value = '''25.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -9.99%  '
$ws.Range("D15").Value = '2.891.12'
$ws.Range("E15").Value = '  -6.90%  '
$ws.Range("D16").Value = '''0.0000164'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -11.20%  '
$ws.Range("D17").Value = '60.961.04'
$ws.Range("E17").Value = '  -4.31%  '
$ws.Range("D18").Value = '2.454.93'
$ws.Range("E18").Value = '  -6.04%  '
$ws.Range("D19").Value = '''10.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -9.71%  '
$ws.Range("D20").Value = '''6.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -9.75%  '
$ws.Range("D21").Value = '''4.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.53%  '
$ws.Range("D22").Value = '''317.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.69%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").Value = '''1.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.05%  '
$ws.Range("D25").Value = '''63.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.79%  '
$ws.Range("D26").Value = '2.580.17'
$ws.Range("E26").Value = '  -6.61%  '
$ws.Range("D27").Value = '0.0₃0953'
$ws.Range("E27").Value = '  -15.64%  '
# Row 28: coin reordering (name/link swap) + refreshed price/volume
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '''1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.04%  '
# Row 29: coin reordering (name/link swap) + refreshed price/volume
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").Value = '''532.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.96%  '
$ws.Range("D30").Value = '''1.44'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -12.71%  '
$ws.Range("D31").Value = '''8.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.90%  '
$ws.Range("D32").Value = '''7.53'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.94%  '
$ws.Range("D33").Value = '''0.146'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.63%  '
$ws.Range("D34").Value = '''1.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.03%  '
$ws.Range("D35").Value = '''1.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.73%  '
$ws.Range("D36").Value = '''5.74'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -13.74%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '''4.76'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -12.90%  '
$ws.Range("D39").Value = '''0.375'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.08%  '
$ws.Range("D40").Value = '''18.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.21%  '
# Row 41: coin reordering (name/link swap) + refreshed price/volume
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").Value = '''142.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.76%  '
# Row 42: coin reordering (name/link swap) + refreshed price/volume
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''1.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.84%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").Value = '''40.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.03%  '
$ws.Range("D45").Value = '''2.26'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -11.84%  '
$ws.Range("D46").Value = '''143.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -12.06%  '
$ws.Range("D47").Value = '''3.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.19%  '
$ws.Range("D48").Value = '''21.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -12.03%  '
$ws.Range("D49").Value = '''0.0529'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.92%  '
$ws.Range("D50").Value = '''0.585'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.56%  '
$ws.Range("D51").Value = '''0.0926'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.47%  '
